# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Re-sort the "Estado de Cuenta" table (rows 16-32) so entries are grouped
# by Periodo Mora (column E) instead of by worker (columns C/D), and add
# the new ALBEIRO JOSE COMAS MARTINEZ / 1902 record that was missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico
$rows = @(
    @("CC","9144662","DIOMEDES DE JESUS TORRES HERNANDEZ","1811",89344,2233590),
    @("CC","1096193978","JOHAN ARLEY GARCIA ESPARZA","1811",85085,2127120),
    @("CC","8867001","ALFONSO LUIS GALARCIO FURNIELES","1811",81031,2025780),
    @("CC","9144662","DIOMEDES DE JESUS TORRES HERNANDEZ","1812",89344,2233590),
    @("CC","1096193978","JOHAN ARLEY GARCIA ESPARZA","1812",85085,2127120),
    @("CC","8867001","ALFONSO LUIS GALARCIO FURNIELES","1812",81031,2025780),
    @("CC","9144662","DIOMEDES DE JESUS TORRES HERNANDEZ","1901",89344,2233590),
    @("CC","1096193978","JOHAN ARLEY GARCIA ESPARZA","1901",85085,2127120),
    @("CC","8867001","ALFONSO LUIS GALARCIO FURNIELES","1901",81031,2025780),
    @("CC","9144662","DIOMEDES DE JESUS TORRES HERNANDEZ","1902",89344,2233590),
    @("CC","9023326","ALBEIRO JOSE COMAS MARTINEZ","1902",81031,2025780),
    @("CC","1096193978","JOHAN ARLEY GARCIA ESPARZA","1902",85085,2127120),
    @("CC","8867001","ALFONSO LUIS GALARCIO FURNIELES","1902",81031,2025780),
    @("CC","9144662","DIOMEDES DE JESUS TORRES HERNANDEZ","1903",74453,2233590),
    @("CC","9023326","ALBEIRO JOSE COMAS MARTINEZ","1903",67526,2025780),
    @("CC","1096193978","JOHAN ARLEY GARCIA ESPARZA","1903",70904,2127120),
    @("CC","8867001","ALFONSO LUIS GALARCIO FURNIELES","1903",67526,2025780)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
}
